$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match-data (columns F:V) between row pairs that were reordered ---
# The "Indice" (A) and "data_partida" (E) columns stay put per physical row;
# everything else (home/away teams, scores, odds, timestamps, url) swaps.

function Swap-Rows($r1, $r2) {
    $range1 = $ws.Range("F" + $r1 + ":V" + $r1)
    $range2 = $ws.Range("F" + $r2 + ":V" + $r2)
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

Swap-Rows 96 97
Swap-Rows 110 111
Swap-Rows 133 134

# --- Append new row 172 (new match result scraped into the sheet) ---
# Copy row 171 first so text/number styles & types (e.g. the "2023" text
# in column D) are preserved, then overwrite with the new match's data.
$ws.Range("A171:V171").Copy($ws.Range("A172:V172"))

$ws.Range("A172").Value = 171
$ws.Range("E172").Value = 45241.83333333334
$ws.Range("F172").Value = "Argentinos Jrs"
$ws.Range("G172").Value = 1
$ws.Range("H172").Value = "Velez Sarsfield"
$ws.Range("I172").Value = 1
$ws.Range("J172").Value = 2.04
$ws.Range("K172").Value = "07/11/2023 06:42"
$ws.Range("L172").Value = 2.2
$ws.Range("M172").Value = "11/11/2023 19:59"
$ws.Range("N172").Value = 3.16
$ws.Range("O172").Value = "07/11/2023 06:42"
$ws.Range("P172").Value = 2.94
$ws.Range("Q172").Value = "11/11/2023 19:59"
$ws.Range("R172").Value = 4.21
$ws.Range("S172").Value = "07/11/2023 06:42"
$ws.Range("T172").Value = 4.19
$ws.Range("U172").Value = "11/11/2023 19:59"
$ws.Range("V172").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/argentinos-jrs-velez-sarsfield/0t0VKDtD/"
